# Update "2025" row (row 8) metrics in the recorrência anual sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1068
$ws.Range("D8").Value = 175
$ws.Range("E8").Value = 893
$ws.Range("F8").Value = 7.178014766201805
$ws.Range("G8").Value = 83.61423220973782
$ws.Range("H8").Value = 16.38576779026217
